$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that currently follows the
#    "Play 9 Lions Slot Game for Free - Review" Heading1 paragraph.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph "Play 9 Lions Slot Game for Free - Review"
#    right before the final "Prompt: ..." paragraph (i.e. just before
#    <w:sectPr>). We build the new paragraph via InsertXML so that it gets
#    exactly the formatting we want (bold only, no inherited italics),
#    matching the structure used elsewhere in the document
#    (<w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>...</w:t></w:r>).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertStart = $lastPara.Range.Start
$insertPoint = $d.Range($insertStart, $insertStart)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 9 Lions Slot Game for Free - Review</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xml)

# InsertXML above leaves behind an extra empty paragraph (used only to force
# a paragraph break) right after our new paragraph; remove it.
$newCount = $d.Paragraphs.Count
$strayPara = $d.Paragraphs.Item($newCount - 1)
$strayPara.Range.Delete()

# 3. Replace the old "Prompt: ..." image-generation text with the new meta
#    description text, keeping the existing italic run formatting intact.
$oldText = 'Prompt: Design a feature image for the online slot game "9 Lions". The image should be in a cartoon style, and should feature a happy Maya warrior with glasses. The Maya warrior should be surrounded by the 9 lions and other Chinese-themed symbols. The background should portray a mystical mountain with the temple dedicated to the 9 lions, with some dragon icons and Chinese lanterns. The overall theme of the image should be bright and colorful, capturing the essence of the beautiful graphics in the game.'
$newText = 'Discover the enchanting graphics and bonus features of the 9 Lions slot game that make it mesmerizing to play. Try it for free now!'
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
